$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete trailing rows (old rows 9 and 10)
$ws.Range("A9:H10").Delete()

# New data (rows 2-8), replacing the atypical-sales dataset with the
# refreshed cronjob output.
# Columns: A=index, B=Dia, C=quantidade_atipica, D=cliente, E=id_produto,
#          F=produto, G=estoque_atualizado, H=critico
$data = @(
    @(0, "2025-04-14", 60,  "AMAZONIA REFEICOES E SERVICOS LTDA",  "000103", "AVENTAL PVC FORRADO PLUS 1,20X0,65 BRANCO C.A. 28303 BRASCAMP", 64,   $false),
    @(3, "2025-04-17", 84,  "METALURGICA SATO DA AMAZONIA LTDA",   "000032", "LIMPADOR VEJA MULTIUSO GOLD 500ML",                              740,  $false),
    @(5, "2025-04-22", 24,  "BRAGA MOTOS LTDA",                    "000015", "PANO MULTIUSO ROLO 28X300 M AZUL TALGE",                         0,    $false),
    @(4, "2025-04-23", 480, "MANJAR SERVICOS GERAIS SA",           "000184", "SACO PLAST 1 KG 16X30 ALTA DENSIDADE C/100 UND",                 -290, $false),
    @(6, "2025-04-23", 400, "V V REFEICOES LTDA",                  "000029", "ESPONJA MULTIUSO JEITOSA",                                       1184, $false),
    @(1, "2025-04-24", 150, "JURUA ESTALEIROS E NAVEGACAO LTDA",   "000088", "VASSOURA PIACAVA 20 FUROS",                                      -16,  $false),
    @(2, "2025-04-24", 300, "MUSASHI DA AMAZONIA LTDA",            "000842", "SACO DE LIXO 200L COMUM PACOTINHO C/5 UND",                      31,   $false)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]

    # B (Dia) and E (id_produto) hold strings that look like dates /
    # zero-padded numbers; force text formatting so Excel doesn't
    # reinterpret them, then strip the leftover number-format so the
    # cell keeps the workbook's default (unstyled) look.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 5).ClearFormats()

    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $row++
}
